$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 491, shifting the existing data (old rows 491:573)
# down to 496:578.
$ws.Rows.Item(491).Resize(5).Insert()

# Constant columns shared by every data row in this sheet.
$constA = 2
$constB = "Comercializadora del Agro de Limarí"
$constC = "Coquimbo"
$constE = 4
$constF = "Fruta"
$constG = 100102
$constH = "Cítricos"
$constI = 100102005
$constJ = "Naranja"
$constQ = "`$/bins (400 kilos)"
$constR = "Provincia de Limarí"
$constT = 400

# New weekly rows (D, K, L, M, N, O, P, S) for rows 491-495.
$newRows = @(
    @{ Row = 491; D = 45218; K = "Cara cara";  L = "Primera"; M = 20; N = 160000; O = 170000; P = 165000; S = 412 },
    @{ Row = 492; D = 45218; K = "Lane Late";  L = "Primera"; M = 20; N = 160000; O = 170000; P = 165000; S = 412 },
    @{ Row = 493; D = 45218; K = "Lane Late";  L = "Segunda"; M = 16; N = 120000; O = 130000; P = 125000; S = 312 },
    @{ Row = 494; D = 45218; K = "Navel Late"; L = "Primera"; M = 16; N = 160000; O = 170000; P = 165000; S = 412 },
    @{ Row = 495; D = 45218; K = "Navel Late"; L = "Segunda"; M = 10; N = 120000; O = 130000; P = 125000; S = 312 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $constI
    $ws.Cells.Item($row, 10).Value = $constJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $constT
}
